$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 216; existing rows 216..307 shift down to 217..308.
$ws.Rows.Item(216).Insert()

# Fill the new row 216 with its data. Columns A,B,C,E,F,G,H,N,Q,R carry the same
# constant values used throughout this dataset; D,I,J,K,L,M,O,P are the new record's
# own values.
$ws.Cells.Item(216, 1).Value = 7
$ws.Cells.Item(216, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(216, 3).Value = "Ñuble"
$ws.Cells.Item(216, 4).Value = 44559
$ws.Cells.Item(216, 5).Value = 16
$ws.Cells.Item(216, 6).Value = 100114014
$ws.Cells.Item(216, 7).Value = "Betarraga"
$ws.Cells.Item(216, 8).Value = "Sin especificar"
$ws.Cells.Item(216, 9).Value = "Primera"
$ws.Cells.Item(216, 10).Value = 800
$ws.Cells.Item(216, 11).Value = 500
$ws.Cells.Item(216, 12).Value = 600
$ws.Cells.Item(216, 13).Value = 550
$ws.Cells.Item(216, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(216, 15).Value = "Región del Maule"
$ws.Cells.Item(216, 16).Value = 110
$ws.Cells.Item(216, 17).Value = 5
$ws.Cells.Item(216, 18).Value = "Hortaliza"
